$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = "e11abc11-d81e-48f3-a5af-a03605f53f72.md"
$wsOverview.Range("A3").Value = "ffff8d28186a-9f8f-4d6e-af74-3eea744c3452.md"

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("A2").Value = "e11abc11-d81e-48f3-a5af-a03605f53f72.md"
$wsZh.Range("D2").Value = "e11abc11-d81e-48f3-a5af-a03605f53f72.e5c2defb3d53fad5e7a5b8091116ba601f6445fc.zh-cn.xlf"
$wsZh.Range("E2").Value = "2016-03-21 12:56:32"
$wsZh.Range("F2").Value = "e11abc11-d81e-48f3-a5af-a03605f53f72.md"
$wsZh.Range("G2").Value = "e11abc11-d81e-48f3-a5af-a03605f53f72.e5c2defb3d53fad5e7a5b8091116ba601f6445fc.zh-cn.xlf"
$wsZh.Range("H2").Value = "2016-03-21 12:56:54"

$wsZh.Range("A3").Value = "ffff8d28186a-9f8f-4d6e-af74-3eea744c3452.md"
$wsZh.Range("D3").Value = "e11abc11-d81e-48f3-a5af-a03605f53f72.e5c2defb3d53fad5e7a5b8091116ba601f6445fc.zh-cn.xlf"
$wsZh.Range("E3").Value = "2016-03-21 12:56:32"
$wsZh.Range("F3").Value = "ffff8d28186a-9f8f-4d6e-af74-3eea744c3452.md"
$wsZh.Range("G3").Value = "e11abc11-d81e-48f3-a5af-a03605f53f72.e5c2defb3d53fad5e7a5b8091116ba601f6445fc.zh-cn.xlf"
$wsZh.Range("H3").Value = "2016-03-21 12:56:54"

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("A2").Value = "e11abc11-d81e-48f3-a5af-a03605f53f72.md"
$wsDe.Range("D2").Value = "e11abc11-d81e-48f3-a5af-a03605f53f72.e5c2defb3d53fad5e7a5b8091116ba601f6445fc.de-de.xlf"
$wsDe.Range("E2").Value = "2016-03-21 12:56:36"
$wsDe.Range("F2").Value = "e11abc11-d81e-48f3-a5af-a03605f53f72.md"
$wsDe.Range("G2").Value = "e11abc11-d81e-48f3-a5af-a03605f53f72.e5c2defb3d53fad5e7a5b8091116ba601f6445fc.de-de.xlf"
$wsDe.Range("H2").Value = "2016-03-21 12:57:00"

$wsDe.Range("A3").Value = "ffff8d28186a-9f8f-4d6e-af74-3eea744c3452.md"
$wsDe.Range("D3").Value = "e11abc11-d81e-48f3-a5af-a03605f53f72.e5c2defb3d53fad5e7a5b8091116ba601f6445fc.de-de.xlf"
$wsDe.Range("E3").Value = "2016-03-21 12:56:36"
$wsDe.Range("F3").Value = "ffff8d28186a-9f8f-4d6e-af74-3eea744c3452.md"
$wsDe.Range("G3").Value = "e11abc11-d81e-48f3-a5af-a03605f53f72.e5c2defb3d53fad5e7a5b8091116ba601f6445fc.de-de.xlf"
$wsDe.Range("H3").Value = "2016-03-21 12:57:00"
